$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells we update so numeric-looking strings (with
# trailing zeros, thousand-dot separators, etc.) are preserved exactly as text,
# matching the source workbook where every data cell is stored as inline text.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("B11").NumberFormat = "@"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("B12").NumberFormat = "@"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("E51").NumberFormat = "@"

$ws.Range("D2").Value = "47.751.74"
$ws.Range("E2").Value = "  +0.91%  "
$ws.Range("D3").Value = "2.508.11"
$ws.Range("E3").Value = "  +0.26%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "322.80"
$ws.Range("E5").Value = "  -0.37%  "
$ws.Range("D6").Value = "108.06"
$ws.Range("E6").Value = "  -1.03%  "
$ws.Range("E7").Value = "  -0.52%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +3.14%  "
$ws.Range("D10").Value = "40.20"
$ws.Range("E10").Value = "  +3.41%  "
$ws.Range("B11").Value = "Chainlink"
$ws.Range("C11").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D11").Value = "19.67"
$ws.Range("E11").Value = "  +6.13%  "
$ws.Range("B12").Value = "Dogecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D12").Value = "0.0814"
$ws.Range("E12").Value = "  -0.32%  "
$ws.Range("E13").Value = "  +0.72%  "
$ws.Range("E14").Value = "  -0.54%  "
$ws.Range("D15").Value = "2.898.23"
$ws.Range("E15").Value = "  +0.11%  "
$ws.Range("D16").Value = "2.509.58"
$ws.Range("E16").Value = "  +0.40%  "
$ws.Range("E17").Value = "  -1.04%  "
$ws.Range("D18").Value = "47.649.02"
$ws.Range("E18").Value = "  +0.79%  "
$ws.Range("D19").Value = "13.35"
$ws.Range("E19").Value = "  +2.80%  "
$ws.Range("D21").Value = "0.0₃0941"
$ws.Range("E21").Value = "  -0.70%  "
$ws.Range("D22").Value = "2.76"
$ws.Range("E22").Value = "  +7.71%  "
$ws.Range("D23").Value = "70.93"
$ws.Range("D24").Value = "247.36"
$ws.Range("E24").Value = "  -1.10%  "
$ws.Range("D25").Value = "2.58"
$ws.Range("E25").Value = "  -1.08%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").Value = "25.76"
$ws.Range("E27").Value = "  -1.47%  "
$ws.Range("D28").Value = "10.22"
$ws.Range("E28").Value = "  +1.95%  "
$ws.Range("E29").Value = "  +4.22%  "
$ws.Range("D30").Value = "34.82"
$ws.Range("E30").Value = "  -2.57%  "
$ws.Range("D32").Value = "49.85"
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("D33").Value = "20.07"
$ws.Range("E33").Value = "  +0.37%  "
$ws.Range("E34").Value = "  -1.58%  "
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("D36").Value = "0.0784"
$ws.Range("E36").Value = "  -1.21%  "
$ws.Range("E37").Value = "  -1.10%  "
$ws.Range("D38").Value = "4.69"
$ws.Range("E38").Value = "  -1.42%  "
$ws.Range("E39").Value = "  -0.60%  "
$ws.Range("E40").Value = "  -0.23%  "
$ws.Range("D41").Value = "22.19"
$ws.Range("E41").Value = "  +3.40%  "
$ws.Range("E42").Value = "  -2.38%  "
$ws.Range("D43").Value = "118.93"
$ws.Range("E43").Value = "  -2.59%  "
$ws.Range("E44").Value = "  -0.60%  "
$ws.Range("D45").Value = "2.002.36"
$ws.Range("E45").Value = "  +0.64%  "
$ws.Range("D46").Value = "3.10"
$ws.Range("E46").Value = "  +1.13%  "
$ws.Range("E47").Value = "  -3.11%  "
$ws.Range("E48").Value = "  +0.79%  "
$ws.Range("D49").Value = "9.09"
$ws.Range("E49").Value = "  +0.33%  "
$ws.Range("E50").Value = "  -3.34%  "
$ws.Range("D51").Value = "56.57"
$ws.Range("E51").Value = "  +2.04%  "
